# "removing thread sleeps for search testcases"
# Re-run results for the B-suite "Test Cases" sheet: reset every test case's
# Results cell (column E, rows 2-89) to "SKIP", except TestCase_B10 (row 11 -
# "Verify that filtering is retained when user navigates back to search
# results page from record view page"), which now comes back "FAIL".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Reset the whole Results column to SKIP (clears the old FAIL in E2, the old
# PASS in E84, and fills in the previously-blank E85:E89).
$ws.Range("E2:E89").Value = "SKIP"

# TestCase_B10 now fails.
$ws.Range("E11").Value = "FAIL"

# Select the Results column that was just updated (this also scrolls the
# view back to the top, clearing the old topLeftCell="A87" scroll state).
$ws.Range("D2:D89").Select() | Out-Null
